$wb = $excel.ActiveWorkbook

$oldBuilt = "built on January 30 2026 16.19.47 EST"
$newBuilt = "built on February 02 2026 12.49.33 EST"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2")
$a2.Value = $a2.Text.Replace($oldBuilt, $newBuilt)

$a6 = $wsAbout.Range("A6")
$a6.Value = $a6.Text.Replace($oldBuilt, $newBuilt)

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 11; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = build_version
    $cell.Value = $cell.Text.Replace($oldBuilt, $newBuilt)
}
